$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.3164144733447235
$ws.Range("D2").Value = 0.7546706271312429

$ws.Range("C3").Value = -0.1114297882184288
$ws.Range("D3").Value = 0.9122859279064721

$ws.Range("C4").Value = -0.7345520746105022
$ws.Range("D4").Value = 0.4703700354543763

$ws.Range("C5").Value = -0.2961244517092343
$ws.Range("D5").Value = 0.7699123115987083

$ws.Range("C6").Value = 0.2310234292671358
$ws.Range("D6").Value = 0.8194316042476113

$ws.Range("C7").Value = -0.7353370720682859
$ws.Range("D7").Value = 0.4699015092556484

$ws.Range("C8").Value = 0.0002650417432753094
$ws.Range("D8").Value = 0.9997909159681173

$ws.Range("C9").Value = -0.8463954232197843
$ws.Range("D9").Value = 0.4064450057275668

$ws.Range("C10").Value = -0.1865365646018382
$ws.Range("D10").Value = 0.8537334251055884

$ws.Range("C11").Value = 0.3877730883822951
$ws.Range("D11").Value = 0.7019104481126832
